$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "The applicant's projects extensively cover web development using ReactJS and includes experience in Django, Python, and JavaScript, which align with the required skills. However, the lack of experience with MongoDB and NodeJS, as required by the job, may have led to a slightly lower score."
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = "The applicant shows strong adaptability, communication, teamwork, and problem-solving skills, as evidenced by the answers provided. Additionally, the willingness to work in Japan and openness to new cultures demonstrate a high level of adaptability and cultural awareness."

# Row 3
$ws.Range("C3").Value = "The applicant has a good score because they have experience with NodeJS, ExpressJS, ReactJS, and Web Development, which align with the company's requirements. However, the lack of experience with MongoDB and JavaScript could be a drawback based on the job description."
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = "The applicant seems adaptable, willing to learn new skills, and is enthusiastic about being part of a new work culture. They also value teamwork and are open to working in Japan, making them a good fit for the company."

# Row 4
$ws.Range("C4").Value = "The applicant has experience in Python, PyTorch, Tensorflow, and ReactJS, which are relevant to the job's requirements. However, the lack of experience in MongoDB, ExpressJS, and NodeJS could have affected the score."
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = "The applicant has displayed strong adaptability and willingness to work in Japan. The responses indicate good communication, teamwork, and problem-solving skills, aligning with the company's soft skill requirements."
